$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roboflow Annotation Report 7/17/2025 - append a new data row to Table1
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Inherit formatting (number format / borders / fill) from the row above,
# the way Excel extends a table's last-row styling onto a newly added row.
$ws.Range("D59:J59").Copy($ws.Range("D60:J60"))
$ws.Rows.Item(60).RowHeight = 15.6

# New day's figures
$ws.Range("D60").Value = "17/7/2025"
$ws.Range("E60").Value = 406
$ws.Range("F60").Value = 924
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 1012
$ws.Range("J60").Value = "N/A"

# Restore the view to where the author left off working
$win = $excel.ActiveWindow
$ws.Range("G64").Select() | Out-Null
$win.ScrollRow = 45
$win.ScrollColumn = 2
